# New crime data collected - weekly CompStat 68th Precinct update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: Volume/Number and report week date range (rich-text cells)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/1/2025  Through  9/7/2025"

# ---------------------------------------------------------------------------
# Helper values used to re-point "placeholder" cells (0 / ***.*) to the
# shared strings already used elsewhere on the sheet, and to restore the
# normal numeric style (s=14) when a placeholder becomes a real number.
# ---------------------------------------------------------------------------

# Row 15 - Rape
$ws.Range("F15").Value = 1
$ws.Range("N15").Value = -47.619047619047

# Row 16 - Robbery  (C/D become "0" placeholders, E becomes "***.*" placeholder)
$ws.Range("C16").Value = "'0"
$ws.Range("D16").Value = "'0"
$ws.Range("E16").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("C16:D16").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("M16").Value = -36.231884057971
$ws.Range("N16").Value = -88.421052631578

# Row 17 - Fel. Assault
$ws.Range("D17").Value = 2
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -53.846153846153
$ws.Range("I17").Value = 101
$ws.Range("J17").Value = 85
$ws.Range("K17").Value = 18.823529411764
$ws.Range("L17").Value = 3.061224489795
$ws.Range("M17").Value = 42.253521126760
$ws.Range("N17").Value = -49.751243781094

# Row 18 - Burglary
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -63.636363636363
$ws.Range("I18").Value = 52
$ws.Range("J18").Value = 62
$ws.Range("K18").Value = -16.129032258064
$ws.Range("L18").Value = -24.637681159420
$ws.Range("M18").Value = -69.411764705882
$ws.Range("N18").Value = -93.157894736842

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 150
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = 75
$ws.Range("I19").Value = 255
$ws.Range("J19").Value = 259
$ws.Range("K19").Value = -1.544401544401
$ws.Range("L19").Value = -15.282392026578
$ws.Range("M19").Value = 19.718309859154
$ws.Range("N19").Value = -17.741935483871

# Row 20 - G.L.A.
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 92
$ws.Range("J20").Value = 142
$ws.Range("K20").Value = -35.211267605633
$ws.Range("L20").Value = 10.843373493975
$ws.Range("M20").Value = -14.814814814814
$ws.Range("N20").Value = -93.087903831705

# Row 21 - TOTAL (bold styles 17/18)
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 30
$ws.Range("F21").Value = 54
$ws.Range("G21").Value = 61
$ws.Range("H21").Value = -11.475409836065
$ws.Range("I21").Value = 555
$ws.Range("J21").Value = 600
$ws.Range("K21").Value = -7.5
$ws.Range("L21").Value = -8.566721581548
$ws.Range("M21").Value = -13.009404388714
$ws.Range("N21").Value = -81.543066178915

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -32.352941176470
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 142
$ws.Range("H24").Value = -37.323943661971
$ws.Range("I24").Value = 710
$ws.Range("J24").Value = 1002
$ws.Range("K24").Value = -29.141716566866
$ws.Range("L24").Value = -32.573599240265
$ws.Range("M24").Value = -10.126582278481

# Row 25 - Retail Theft
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 80
$ws.Range("H25").Value = -56.25
$ws.Range("I25").Value = 326
$ws.Range("J25").Value = 630
$ws.Range("K25").Value = -48.253968253968
$ws.Range("L25").Value = -39.179104477611

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = -14.285714285714
$ws.Range("I26").Value = 275
$ws.Range("J26").Value = 247
$ws.Range("K26").Value = 11.336032388664
$ws.Range("L26").Value = 5.363984674329
$ws.Range("M26").Value = 7.843137254901

# Row 27 - UCR Rape*  (D becomes "0" placeholder, E becomes "***.*" placeholder)
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "'***.*"
$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E23").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = 0

# Row 28 - Other Sex Crimes  (C28 placeholder "0" becomes a real number)
$ws.Range("C28").Value = 2
$ws.Range("D28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = -81.818181818181
$ws.Range("I28").Value = 29
$ws.Range("J28").Value = 41
$ws.Range("K28").Value = -29.268292682926
$ws.Range("L28").Value = 38.095238095238
